# Production Verification Scripts done
# Update the "Date" stamp (column B, row 2) on each of the six bootstrap
# sheets to reflect the latest execution timestamp recorded by the
# Katalon test run.

$wb = $excel.ActiveWorkbook

$stamps = @{
    "CreateUser"    = "Fri Mar 14 07:15:49 IST 2025"
    "FindUser"      = "Fri Mar 14 07:16:42 IST 2025"
    "ModifyUser"    = "Fri Mar 14 07:17:22 IST 2025"
    "ModifyUserPwd" = "Fri Mar 14 07:18:19 IST 2025"
    "AddDeleteRole" = "Fri Mar 14 07:14:12 IST 2025"
    "SearchRole"    = "Fri Mar 14 07:15:04 IST 2025"
}

foreach ($sheetName in $stamps.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B2").Value = $stamps[$sheetName]
}
